# Inserts a new sales-rep row ("LOAIZA TINOCO JUAN PABLO", under the
# "OFICINA-CATAECSA" group) into the alphabetically-sorted employee lists on
# the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. The new row lands right
# before "LOZANO MOLINA TITO JERSON" (between "LINCANGO ..." and
# "LOZANO ..."), pushing every following row down by one. The trailing
# summary row on "VENTAS POR GRUPO" ("N de 326" style counters) is updated to
# reflect the new total headcount of 327.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": columns A..R, new row inserted at row 287.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(287).Insert()
$ws1.Cells.Item(287, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(287, 2).Value = "LOAIZA TINOCO JUAN PABLO"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(287, $c).Value = 0
}

# The last row of the sheet is a "<n> de 326" style summary row that shifted
# from row 328 to row 329 because of the insert above; bump the headcount
# text from 326 to 327 to match the new number of employees.
$totalsRow1 = 329
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item($totalsRow1, $c)
    $cell.Value = $cell.Value2.Replace("de 326", "de 327")
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": columns A..G, new row inserted at row 291.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(291).Insert()
$ws2.Cells.Item(291, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(291, 2).Value = "LOAIZA TINOCO JUAN PABLO"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(291, $c).Value = 0
}
